$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update existing cell values (rows 496, 498, 500, 518) ---
$ws.Cells.Item(496, 4).Value = 14
$ws.Cells.Item(498, 3).Value = 127
$ws.Cells.Item(500, 3).Value = 90
$ws.Cells.Item(518, 3).Value = 10

# --- 2. Copy formatting (styles s=5 / s=6 + row height) from the last existing
#        data row (520) down onto the new rows (521:544) before filling values ---
$ws.Range("A520:F520").Copy()
$ws.Range("A521:F544").PasteSpecial(-4122)

# --- 3. Populate the 24 new data rows ---
$ws.Cells.Item(521, 1).Value = "2024-12-12 23:22:54"
$ws.Cells.Item(521, 2).Value = "021862.NC"
$ws.Cells.Item(521, 3).Value = 100
$ws.Cells.Item(521, 4).Value = 24
$ws.Cells.Item(521, 5).Value = 38
$ws.Cells.Item(521, 6).Value = 1
$ws.Cells.Item(522, 1).Value = "2024-12-17 09:50:32"
$ws.Cells.Item(522, 2).Value = "022048.NC"
$ws.Cells.Item(522, 3).Value = 253
$ws.Cells.Item(522, 4).Value = 133
$ws.Cells.Item(522, 5).Value = 25
$ws.Cells.Item(522, 6).Value = 1
$ws.Cells.Item(523, 1).Value = "2024-12-18 08:46:33"
$ws.Cells.Item(523, 2).Value = "021978.NC"
$ws.Cells.Item(523, 3).Value = 11
$ws.Cells.Item(523, 4).Value = 25
$ws.Cells.Item(523, 5).Value = 12
$ws.Cells.Item(523, 6).Value = 1
$ws.Cells.Item(524, 1).Value = "2024-12-18 11:55:54"
$ws.Cells.Item(524, 2).Value = "021954.NC"
$ws.Cells.Item(524, 3).Value = 6
$ws.Cells.Item(524, 4).Value = 9
$ws.Cells.Item(524, 5).Value = 2
$ws.Cells.Item(524, 6).Value = 1
$ws.Cells.Item(525, 1).Value = "2024-12-18 12:08:47"
$ws.Cells.Item(525, 2).Value = "021955.NC"
$ws.Cells.Item(525, 3).Value = 6
$ws.Cells.Item(525, 4).Value = 9
$ws.Cells.Item(525, 5).Value = 2
$ws.Cells.Item(525, 6).Value = 1
$ws.Cells.Item(526, 1).Value = "2024-12-18 12:34:39"
$ws.Cells.Item(526, 2).Value = "021874.NC"
$ws.Cells.Item(526, 3).Value = 92
$ws.Cells.Item(526, 4).Value = 363
$ws.Cells.Item(526, 5).Value = 6
$ws.Cells.Item(526, 6).Value = 1
$ws.Cells.Item(527, 1).Value = "2024-12-18 13:39:05"
$ws.Cells.Item(527, 2).Value = "022061.NC"
$ws.Cells.Item(527, 3).Value = 10
$ws.Cells.Item(527, 4).Value = 47
$ws.Cells.Item(527, 5).Value = 6
$ws.Cells.Item(527, 6).Value = 1
$ws.Cells.Item(528, 1).Value = "2024-12-18 15:39:58"
$ws.Cells.Item(528, 2).Value = "022019.NC"
$ws.Cells.Item(528, 3).Value = 88
$ws.Cells.Item(528, 4).Value = 420
$ws.Cells.Item(528, 5).Value = 10
$ws.Cells.Item(528, 6).Value = 1
$ws.Cells.Item(529, 1).Value = "2024-12-19 13:10:56"
$ws.Cells.Item(529, 2).Value = "022020.NC"
$ws.Cells.Item(529, 3).Value = 16
$ws.Cells.Item(529, 4).Value = 56
$ws.Cells.Item(529, 5).Value = 10
$ws.Cells.Item(529, 6).Value = 1
$ws.Cells.Item(530, 1).Value = "2024-12-19 16:50:33"
$ws.Cells.Item(530, 2).Value = "022203.NC"
$ws.Cells.Item(530, 3).Value = 2
$ws.Cells.Item(530, 4).Value = 10
$ws.Cells.Item(530, 5).Value = 6
$ws.Cells.Item(530, 6).Value = 1
$ws.Cells.Item(531, 1).Value = "2024-12-20 10:43:11"
$ws.Cells.Item(531, 2).Value = "022087.NC"
$ws.Cells.Item(531, 3).Value = 3
$ws.Cells.Item(531, 4).Value = 17
$ws.Cells.Item(531, 5).Value = 2
$ws.Cells.Item(531, 6).Value = 1
$ws.Cells.Item(532, 1).Value = "2024-12-20 11:47:17"
$ws.Cells.Item(532, 2).Value = "022110.NC"
$ws.Cells.Item(532, 3).Value = 2
$ws.Cells.Item(532, 4).Value = 7
$ws.Cells.Item(532, 5).Value = 2
$ws.Cells.Item(532, 6).Value = 1
$ws.Cells.Item(533, 1).Value = "2024-12-20 12:10:03"
$ws.Cells.Item(533, 2).Value = "022123.NC"
$ws.Cells.Item(533, 3).Value = 1
$ws.Cells.Item(533, 4).Value = 6
$ws.Cells.Item(533, 5).Value = 2
$ws.Cells.Item(533, 6).Value = 1
$ws.Cells.Item(534, 1).Value = "2024-12-20 12:26:12"
$ws.Cells.Item(534, 2).Value = "022133.NC"
$ws.Cells.Item(534, 3).Value = 1
$ws.Cells.Item(534, 4).Value = 11
$ws.Cells.Item(534, 5).Value = 2
$ws.Cells.Item(534, 6).Value = 1
$ws.Cells.Item(535, 1).Value = "2024-12-20 13:00:44"
$ws.Cells.Item(535, 2).Value = "022095.NC"
$ws.Cells.Item(535, 3).Value = 2
$ws.Cells.Item(535, 4).Value = 17
$ws.Cells.Item(535, 5).Value = 2
$ws.Cells.Item(535, 6).Value = 1
$ws.Cells.Item(536, 1).Value = "2024-12-20 13:17:21"
$ws.Cells.Item(536, 2).Value = "022086.NC"
$ws.Cells.Item(536, 3).Value = 2
$ws.Cells.Item(536, 4).Value = 17
$ws.Cells.Item(536, 5).Value = 2
$ws.Cells.Item(536, 6).Value = 1
$ws.Cells.Item(537, 1).Value = "2024-12-20 13:28:30"
$ws.Cells.Item(537, 2).Value = "022094.NC"
$ws.Cells.Item(537, 3).Value = 3
$ws.Cells.Item(537, 4).Value = 17
$ws.Cells.Item(537, 5).Value = 2
$ws.Cells.Item(537, 6).Value = 1
$ws.Cells.Item(538, 1).Value = "2024-12-20 13:49:34"
$ws.Cells.Item(538, 2).Value = "022093.NC"
$ws.Cells.Item(538, 3).Value = 1
$ws.Cells.Item(538, 4).Value = 17
$ws.Cells.Item(538, 5).Value = 2
$ws.Cells.Item(538, 6).Value = 1
$ws.Cells.Item(539, 1).Value = "2024-12-20 14:58:26"
$ws.Cells.Item(539, 2).Value = "022089.NC"
$ws.Cells.Item(539, 3).Value = 3
$ws.Cells.Item(539, 4).Value = 17
$ws.Cells.Item(539, 5).Value = 2
$ws.Cells.Item(539, 6).Value = 1
$ws.Cells.Item(540, 1).Value = "2024-12-20 15:12:20"
$ws.Cells.Item(540, 2).Value = "022092.NC"
$ws.Cells.Item(540, 3).Value = 4
$ws.Cells.Item(540, 4).Value = 17
$ws.Cells.Item(540, 5).Value = 2
$ws.Cells.Item(540, 6).Value = 1
$ws.Cells.Item(541, 1).Value = "2024-12-20 15:51:20"
$ws.Cells.Item(541, 2).Value = "022090.NC"
$ws.Cells.Item(541, 3).Value = 3
$ws.Cells.Item(541, 4).Value = 17
$ws.Cells.Item(541, 5).Value = 2
$ws.Cells.Item(541, 6).Value = 1
$ws.Cells.Item(542, 1).Value = "2024-12-23 09:20:09"
$ws.Cells.Item(542, 2).Value = "022096.NC"
$ws.Cells.Item(542, 3).Value = 4
$ws.Cells.Item(542, 4).Value = 18
$ws.Cells.Item(542, 5).Value = 2
$ws.Cells.Item(542, 6).Value = 1
$ws.Cells.Item(543, 1).Value = "2024-12-23 10:26:46"
$ws.Cells.Item(543, 2).Value = "022135.NC"
$ws.Cells.Item(543, 3).Value = 1
$ws.Cells.Item(543, 4).Value = 9
$ws.Cells.Item(543, 5).Value = 2
$ws.Cells.Item(543, 6).Value = 1
$ws.Cells.Item(544, 1).Value = "2024-12-23 10:35:17"
$ws.Cells.Item(544, 2).Value = "022136.NC"
$ws.Cells.Item(544, 3).Value = 2
$ws.Cells.Item(544, 4).Value = 15
$ws.Cells.Item(544, 5).Value = 2
$ws.Cells.Item(544, 6).Value = 1

# Row heights for the new rows (13.55pt, matching the rest of the sheet)
$ws.Range("A521:F544").RowHeight = 13.55

# --- 4. Column width changes: column A narrows, columns B:F take on the shared
#        default data width (best achievable through the ColumnWidth property,
#        which snaps to Excel pixel-width quantisation) ---
$ws.Range("A1").EntireColumn.ColumnWidth = 18
$ws.Range("B1:F1").EntireColumn.ColumnWidth = 8
